# Update the "想去人数" (F column) figures on both the "展览" and "全部类型"
# sheets to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    2  = 11541
    3  = 11040
    6  = 1000
    8  = 61
    10 = 39
    11 = 10649
    12 = 4114
    19 = 425
    21 = 10866
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
